# Update the per-class statistics table (rows 2-11, columns A:C) with the
# newly computed selectivity / sensitivity / accuracy values.
#
# These values are stored in the workbook as text (shared strings), not as
# numbers, so we temporarily force the Text number format before writing
# them and then clear the format again so the cells end up with the default
# style (matching the original workbook, where no cell carries an explicit
# style index).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dataRange = $ws.Range("A2:C11")
$dataRange.NumberFormat = "@"

$values = @(
    @("0.92", "0.75", "0.83"),
    @("0.76", "1.00", "0.86"),
    @("0.84", "1.00", "0.91"),
    @("0.80", "1.00", "0.89"),
    @("0.94", "1.00", "0.97"),
    @("0.88", "0.88", "0.88"),
    @("0.75", "0.94", "0.83"),
    @("0.93", "0.88", "0.90"),
    @("0.93", "0.81", "0.87"),
    @("0.92", "0.75", "0.83")
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = 2 + $i
    $rowValues = $values[$i]
    $ws.Range("A$row").Value = $rowValues[0]
    $ws.Range("B$row").Value = $rowValues[1]
    $ws.Range("C$row").Value = $rowValues[2]
}

$dataRange.ClearFormats()
